$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Slide 5: change the table's style (tableStyleId) to a different
#    built-in table style.
# ------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tableShape = $s5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{C673E6CE-FD24-4D69-AC99-3C0A4AB8FD2C}")

# ------------------------------------------------------------------
# 2) Re-colour the presentation theme (theme1.xml) from the old
#    "Integral / Red Violet" palette to the standard "Office" palette.
#    ThemeColorScheme items map 1:1 onto the clrScheme children:
#      1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#      8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
#    RGB is packed the classic OLE/VBA way: R + G*256 + B*65536
# ------------------------------------------------------------------
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0x00 + (0x00 * 256) + (0x00 * 65536)   # dk1      000000
$colorScheme.Item(2).RGB  = 0xFF + (0xFF * 256) + (0xFF * 65536)   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 0x44 + (0x54 * 256) + (0x6A * 65536)   # dk2      44546A
$colorScheme.Item(4).RGB  = 0xE7 + (0xE6 * 256) + (0xE6 * 65536)   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 0x5B + (0x9B * 256) + (0xD5 * 65536)   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 0xED + (0x7D * 256) + (0x31 * 65536)   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 0xA5 + (0xA5 * 256) + (0xA5 * 65536)   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 0xFF + (0xC0 * 256) + (0x00 * 65536)   # accent4  FFC000
$colorScheme.Item(9).RGB  = 0x44 + (0x72 * 256) + (0xC4 * 65536)   # accent5  4472C4
$colorScheme.Item(10).RGB = 0x70 + (0xAD * 256) + (0x47 * 65536)   # accent6  70AD47
$colorScheme.Item(11).RGB = 0x05 + (0x63 * 256) + (0xC1 * 65536)   # hlink    0563C1
$colorScheme.Item(12).RGB = 0x95 + (0x4F * 256) + (0x72 * 65536)   # folHlink 954F72
